# Auto-generated Excel COM-interop script
# Applies numeric corrections to H:N profit columns across 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2708540.2
$ws.Range("J17").Value = 2947433
$ws.Range("L17").Value = 8842299
$ws.Range("N17").Value = -8842635
$ws.Range("H116").Value = 4303.6113
$ws.Range("I116").Value = 1975.6666
$ws.Range("J116").Value = 5467.5835
$ws.Range("K116").Value = 1975.6666
$ws.Range("L116").Value = 5467.5835
$ws.Range("M116").Value = 1466.3334
$ws.Range("N116").Value = -12351.5835
$ws.Range("H129").Value = 313753.4
$ws.Range("J129").Value = 371794.4
$ws.Range("L129").Value = 1115383.2
$ws.Range("N129").Value = -1125383.2
$ws.Range("H132").Value = 2363.3
$ws.Range("I132").Value = 2779.5806
$ws.Range("J132").Value = 929.44446
$ws.Range("K132").Value = 8338.7418
$ws.Range("L132").Value = 2788.33338
$ws.Range("M132").Value = -5808.7418
$ws.Range("N132").Value = -7848.33338
$ws.Range("H135").Value = 9618122
$ws.Range("I135").Value = 472.2381
$ws.Range("K135").Value = 4250.1429
$ws.Range("M135").Value = -1715.1429
$ws.Range("H137").Value = 1272.25
$ws.Range("I137").Value = 1251.9565
$ws.Range("J137").Value = 1365.6
$ws.Range("K137").Value = 3755.8695
$ws.Range("L137").Value = 4096.799999999999
$ws.Range("M137").Value = -1205.8695
$ws.Range("N137").Value = -9196.799999999999
$ws.Range("H138").Value = 21278760
$ws.Range("I138").Value = 30304184
$ws.Range("J138").Value = 4548.357
$ws.Range("K138").Value = 90912552
$ws.Range("L138").Value = 13645.071
$ws.Range("M138").Value = -90907412
$ws.Range("N138").Value = -23925.071
$ws.Range("H141").Value = 1315.551
$ws.Range("I141").Value = 804.0513
$ws.Range("K141").Value = 2412.1539
$ws.Range("M141").Value = 2767.8461

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3472.6963
$ws.Range("I32").Value = 3136.6812
$ws.Range("J32").Value = 5791.2
$ws.Range("K32").Value = 3136.6812
$ws.Range("L32").Value = 5791.2
$ws.Range("M32").Value = -2849.6812
$ws.Range("N32").Value = -6365.2
$ws.Range("H45").Value = 2685.2964
$ws.Range("I45").Value = 2893.4
$ws.Range("J45").Value = 2425.1667
$ws.Range("K45").Value = 2893.4
$ws.Range("L45").Value = 2425.1667
$ws.Range("M45").Value = -2516.4
$ws.Range("N45").Value = -3179.1667
$ws.Range("H61").Value = 334479.72
$ws.Range("I61").Value = 400925.7
$ws.Range("K61").Value = 400925.7
$ws.Range("M61").Value = -400713.7
$ws.Range("H63").Value = 3005
$ws.Range("I63").Value = 3005
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3005
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2319
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3005
$ws.Range("I66").Value = 3005
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 15025
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11593
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 28573170
$ws.Range("I74").Value = 35716224
$ws.Range("J74").Value = 953.2857
$ws.Range("K74").Value = 35716224
$ws.Range("L74").Value = 953.2857
$ws.Range("M74").Value = -35715350
$ws.Range("N74").Value = -2701.2857
$ws.Range("H77").Value = 28573170
$ws.Range("I77").Value = 35716224
$ws.Range("J77").Value = 953.2857
$ws.Range("K77").Value = 178581120
$ws.Range("L77").Value = 4766.4285
$ws.Range("M77").Value = -178576752
$ws.Range("N77").Value = -13502.4285
$ws.Range("H132").Value = 14906.026
$ws.Range("I132").Value = 1690.1
$ws.Range("J132").Value = 64465.75
$ws.Range("K132").Value = 5070.299999999999
$ws.Range("L132").Value = 193397.25
$ws.Range("M132").Value = -2540.299999999999
$ws.Range("N132").Value = -198457.25
$ws.Range("H136").Value = 334479.72
$ws.Range("I136").Value = 400925.7
$ws.Range("K136").Value = 1202777.1
$ws.Range("M136").Value = -1200227.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 12831.728
$ws.Range("J81").Value = 12831.728
$ws.Range("L81").Value = 12831.728
$ws.Range("N81").Value = -14953.728
$ws.Range("H84").Value = 12831.728
$ws.Range("J84").Value = 12831.728
$ws.Range("L84").Value = 38495.18399999999
$ws.Range("N84").Value = -49103.18399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3377.7708
$ws.Range("I31").Value = 2076.0625
$ws.Range("J31").Value = 5981.1875
$ws.Range("K31").Value = 2076.0625
$ws.Range("L31").Value = 5981.1875
$ws.Range("M31").Value = -1781.0625
$ws.Range("N31").Value = -6571.1875
$ws.Range("H34").Value = 3377.7708
$ws.Range("I34").Value = 2076.0625
$ws.Range("J34").Value = 5981.1875
$ws.Range("K34").Value = 2076.0625
$ws.Range("L34").Value = 5981.1875
$ws.Range("M34").Value = -1874.0625
$ws.Range("N34").Value = -6385.1875
$ws.Range("H58").Value = 11719.106
$ws.Range("I58").Value = 905.54285
$ws.Range("K58").Value = 905.54285
$ws.Range("M58").Value = -702.54285
$ws.Range("H132").Value = 2149.262
$ws.Range("I132").Value = 1393.2368
$ws.Range("K132").Value = 4179.7104
$ws.Range("M132").Value = -1649.7104
$ws.Range("H134").Value = 578.76
$ws.Range("I134").Value = 518
$ws.Range("J134").Value = 1277.5
$ws.Range("K134").Value = 1554
$ws.Range("L134").Value = 3832.5
$ws.Range("M134").Value = 981
$ws.Range("N134").Value = -8902.5
$ws.Range("H136").Value = 11719.106
$ws.Range("I136").Value = 905.54285
$ws.Range("K136").Value = 2716.62855
$ws.Range("M136").Value = -166.6285500000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1579352.4
$ws.Range("I4").Value = 313
$ws.Range("J4").Value = 7500750
$ws.Range("K4").Value = 939
$ws.Range("L4").Value = 22502250
$ws.Range("M4").Value = -827
$ws.Range("N4").Value = -22502474
$ws.Range("H11").Value = 300.2
$ws.Range("I11").Value = 299.75
$ws.Range("J11").Value = 302
$ws.Range("K11").Value = 899.25
$ws.Range("L11").Value = 906
$ws.Range("M11").Value = -759.25
$ws.Range("N11").Value = -1186
$ws.Range("H18").Value = 234.44444
$ws.Range("I18").Value = 201.25
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 603.75
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -434.75
$ws.Range("N18").Value = -1838
$ws.Range("H26").Value = 692
$ws.Range("I26").Value = 230
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 690
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = -402
$ws.Range("N26").Value = -3576
$ws.Range("H41").Value = 708.8
$ws.Range("I41").Value = 514.6667
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 1544.0001
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -1206.0001
$ws.Range("N41").Value = -3676
$ws.Range("H121").Value = 934.53125
$ws.Range("I121").Value = 429.2857
$ws.Range("J121").Value = 1076
$ws.Range("K121").Value = 1287.8571
$ws.Range("L121").Value = 3228
$ws.Range("M121").Value = 22.14289999999983
$ws.Range("N121").Value = -5848
$ws.Range("H131").Value = 783.38
$ws.Range("I131").Value = 540
$ws.Range("J131").Value = 796.18945
$ws.Range("K131").Value = 1620
$ws.Range("L131").Value = 2388.56835
$ws.Range("M131").Value = 3420
$ws.Range("N131").Value = -12468.56835

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3164.8647
$ws.Range("I80").Value = 2523.7693
$ws.Range("J80").Value = 3512.125
$ws.Range("K80").Value = 2523.7693
$ws.Range("L80").Value = 3512.125
$ws.Range("M80").Value = -1525.7693
$ws.Range("N80").Value = -5508.125
$ws.Range("H83").Value = 3164.8647
$ws.Range("I83").Value = 2523.7693
$ws.Range("J83").Value = 3512.125
$ws.Range("K83").Value = 12618.8465
$ws.Range("L83").Value = 17560.625
$ws.Range("M83").Value = -7626.8465
$ws.Range("N83").Value = -27544.625
$ws.Range("H132").Value = 52585.184
$ws.Range("I132").Value = 11070.667
$ws.Range("J132").Value = 102402.6
$ws.Range("K132").Value = 33212.001
$ws.Range("L132").Value = 307207.8
$ws.Range("M132").Value = -30682.001
$ws.Range("N132").Value = -312267.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1709.0193
$ws.Range("I132").Value = 1240.7561
$ws.Range("J132").Value = 3454.3635
$ws.Range("K132").Value = 3722.2683
$ws.Range("L132").Value = 10363.0905
$ws.Range("M132").Value = -1192.2683
$ws.Range("N132").Value = -15423.0905
$ws.Range("H136").Value = 877.3226
$ws.Range("I136").Value = 877.3226
$ws.Range("K136").Value = 2631.9678
$ws.Range("M136").Value = -81.9677999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 27714.5
$ws.Range("J123").Value = 27714.5
$ws.Range("L123").Value = 27714.5
$ws.Range("N123").Value = -37514.5
$ws.Range("H132").Value = 846.4655
$ws.Range("I132").Value = 638.8409
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 1916.5227
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = 613.4773
$ws.Range("N132").Value = -9557
$ws.Range("H136").Value = 15386506
$ws.Range("I136").Value = 22223108
$ws.Range("J136").Value = 4150.5
$ws.Range("K136").Value = 66669324
$ws.Range("L136").Value = 12451.5
$ws.Range("M136").Value = -66666774
$ws.Range("N136").Value = -17551.5
